$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sensor-data category values (German version categories).
$values = @{
  "K2" = "0.4"
  "Q2" = "0.23"
  "H3" = "0.002"
  "Q3" = "0.16"
  "R3" = "0.1"
  "S3" = "62.4"
  "K4" = "0.2"
  "Q4" = "0.05"
  "R4" = "0.14"
  "K5" = "0.2"
  "R5" = "0.18"
  "S5" = "52.1"
  "H6" = "0.0"
  "O6" = "0.0"
  "Q6" = "0.24"
  "R6" = "0.20"
  "S6" = "50.3"
  "V6" = "0.0"
  "H7" = "0.0"
  "K7" = "0.1"
  "O7" = "0.0"
  "Q7" = "0.11"
  "R7" = "0.24"
  "S7" = "48.2"
  "V7" = "0.0"
  "H8" = "0.0"
  "O8" = "0.0"
  "Q8" = "0.14"
  "S8" = "44.3"
  "V8" = "0.0"
  "I9" = "0.001"
  "J9" = "0.08"
  "Q9" = "0.26"
  "R9" = "0.17"
  "S9" = "42.9"
  "I10" = "0.002"
  "J10" = "0.15"
  "Q10" = "0.4"
  "R10" = "0.24"
  "S10" = "41.0"
  "I11" = "0.002"
  "K11" = "0.10"
  "O11" = "0.003"
  "P11" = "0.002"
  "Q11" = "0.23"
  "R11" = "0.28"
  "S11" = "40.0"
  "K12" = "0.09"
  "O12" = "0.004"
  "Q12" = "0.26"
  "R12" = "0.24"
  "S12" = "41.0"
  "K13" = "0.08"
  "O13" = "0.0012"
  "R13" = "0.22"
  "S13" = "42.9"
  "K14" = "0.12"
  "O14" = "0.002"
  "Q14" = "0.10"
  "R14" = "0.2"
  "S14" = "42.1"
  "K15" = "0.3"
  "O15" = "0.001"
  "P15" = "0.001"
  "Q15" = "0.14"
  "R15" = "0.24"
  "S15" = "44.0"
  "H16" = "0.003"
  "K16" = "0.32"
  "O16" = "0.003"
  "P16" = "0.002"
  "Q16" = "0.12"
  "R16" = "0.31"
  "S16" = "48.0"
  "H17" = "0.001"
  "K17" = "0.33"
  "O17" = "0.000"
  "P17" = "0.003"
  "Q17" = "0.14"
  "S17" = "52.0"
  "H18" = "0.001"
  "P18" = "0.001"
  "Q18" = "0.24"
  "S18" = "44.0"
  "H19" = "0.004"
  "K19" = "0.24"
  "O19" = "0.002"
  "P19" = "0.002"
  "Q19" = "0.28"
  "R19" = "0.48"
  "S19" = "54.0"
  "K20" = "0.31"
  "O20" = "0.001"
  "P20" = "0.001"
  "Q20" = "0.34"
  "R20" = "0.52"
  "S20" = "56.0"
  "K21" = "0.24"
  "Q21" = "0.34"
  "R21" = "0.54"
  "S21" = "59.0"
  "K22" = "0.22"
  "Q22" = "0.36"
  "R22" = "0.59"
  "S22" = "62.4"
  "R23" = "0.6"
  "S23" = "68.4"
  "G24" = "52.8"
  "J24" = "0.12"
  "P24" = "0.001"
  "Q24" = "0.5"
  "R24" = "0.7"
  "S24" = "70.4"
  "G25" = "51.4"
  "J25" = "0.08"
  "K25" = "0.16"
  "P25" = "0.005"
  "Q25" = "0.6"
  "R25" = "0.64"
  "S25" = "72.4"
  "G26" = "70.1"
  "J26" = "0.07"
  "P26" = "0.003"
  "Q26" = "0.65"
  "R26" = "0.66"
  "S26" = "74.2"
  "G27" = "80.2"
  "J27" = "0.09"
  "K27" = "0.12"
  "Q27" = "0.64"
  "R27" = "0.67"
  "S27" = "72.4"
  "V27" = "0.001"
  "G28" = "78.1"
  "H28" = "0.006"
  "K28" = "0.3"
  "Q28" = "0.63"
  "R28" = "0.62"
  "S28" = "70.4"
  "V28" = "0.031"
  "G29" = "77.4"
  "H29" = "0.005"
  "K29" = "0.4"
  "Q29" = "0.60"
  "R29" = "0.61"
  "S29" = "72.4"
  "V29" = "0.043"
  "G30" = "74.3"
  "H30" = "0.003"
  "K30" = "0.5"
  "Q30" = "0.54"
  "R30" = "0.54"
  "S30" = "68.3"
  "V30" = "0.054"
  "G31" = "71.1"
  "H31" = "0.001"
  "J31" = "0.2"
  "K31" = "0.6"
  "Q31" = "0.43"
  "R31" = "0.58"
  "S31" = "62.4"
  "V31" = "0.05"
}

foreach ($cellRef in $values.Keys) {
  $cell = $ws.Range($cellRef)
  $cell.NumberFormat = "@"
  $cell.Value = $values[$cellRef]
  $cell.ClearFormats()
}

# Restore view state: scrolled to column G, active selection O17.
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
$ws.Range("O17").Select()